$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (TC6_SearchResults_Typeahead): rename the shared string used in C13 ---
$ws1.Range("C13").Value = "ValidSearchHeader1"

# --- Sheet2 (Testdata): re-touch B6 text so its shared-string slot moves to the
#     end of the table (matches the author's resave ordering) and drop the
#     now-unused search-results hyperlink on B6, keeping only the B2 hyperlink ---
$ws2.Range("B6").Value = "Showing Results for ""sprocket"""

$ws2.Hyperlinks.Delete()
$h = $ws2.Hyperlinks.Add($ws2.Range("B2"), "https://192.168.15.18/storeus")
$h.TextToDisplay = "https://192.168.15.18/storeus"
$ws2.Range("B2").Value = "`$BaseURL"

# --- Restore the active selections recorded in the saved workbook ---
$ws1.Range("D12").Select()
$ws2.Range("G11").Select()
